$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Import Survey Response 1")
$ws.Name = "TEST_IMPORT_SURVEY_RESP_1_test"
